$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 112 (pushes existing rows 112-119 down to 113-120)
$ws.Rows("112:112").Insert()

# Populate the new row 112 with the weekly Jengibre price data
$ws.Range("A112").Value = 8
$ws.Range("B112").Value = "Terminal La Palmera de La Serena"
$ws.Range("C112").Value = "Coquimbo"
$ws.Range("D112").Value = 45021
$ws.Range("E112").Value = 4
$ws.Range("F112").Value = 100114007
$ws.Range("G112").Value = "Jengibre"
$ws.Range("H112").Value = "Sin especificar"
$ws.Range("I112").Value = "Primera"
$ws.Range("J112").Value = 400
$ws.Range("K112").Value = 17000
$ws.Range("L112").Value = 18000
$ws.Range("M112").Value = 17500
$ws.Range("N112").Value = "$/caja 13 kilos"
$ws.Range("O112").Value = "Perú"
$ws.Range("P112").Value = 1346
$ws.Range("Q112").Value = 13
$ws.Range("R112").Value = "Hortaliza"
